$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update the birthday text for row 5 (Howard): 12/05/1987 -> 21/05/2002
$ws1.Range("E5").Value = "21/05/2002"

# Add a new row 6 for a new customer "Billy"
$ws1.Range("A6").Value = "Billy"
$ws1.Range("B6").Value = 23
$ws1.Range("C6").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws1.Range("C6").Value = $true
$ws1.Range("D6").Value = "Beer"
$ws1.Range("E6").NumberFormat = "mm/dd/yy"
$ws1.Range("E6").Formula = "=DATE(1940,4,28)"
$ws1.Range("F6").Value = 1.25

# Make the customers sheet the active tab/sheet, with E7 selected
$ws1.Activate() | Out-Null
$ws1.Range("E7").Select() | Out-Null
